$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 52

# Column A holds a literal date-like text string (not a real date value),
# matching the existing rows (e.g. A10, A12, ...). Force text formatting
# before assignment so Excel doesn't auto-convert the string into a date
# serial number, then restore the default "Normal" style so no stray
# per-cell formatting is left behind (matches the un-styled neighbouring
# rows in the sheet).
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "08/22/2025"
$cellA.Style = "Normal"

$ws.Cells.Item($row, 2).Value = 578.1140000000014
$ws.Cells.Item($row, 3).Value = 0.08648813209851323
$ws.Cells.Item($row, 4).Value = 50
